# Lab1 "dates" workbook: refresh the pressure readings (column B) with a new
# data set, and format that column (centered + wrapped) to match the
# reformatted header/data table described by the commit "Added dates and dates1".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New Pressure values for B2:B49 -----------------------------------
$newVals = @{
    2  = 65
    3  = 70
    4  = 68
    5  = 66
    6  = 69
    7  = 71
    8  = 74
    9  = 90
    10 = 90
    11 = 83
    12 = 81
    13 = 86
    14 = 83
    15 = 80
    16 = 74
    17 = 67
    18 = 64
    19 = 64
    20 = 53
    21 = 61
    22 = 54
    23 = 51
    24 = 64
    25 = 64
    26 = 67
    27 = 68
    28 = 65
    29 = 57
    30 = 59
    31 = 87
    32 = 89
    33 = 76
    34 = 75
    35 = 88
    36 = 90
    37 = 91
    38 = 87
    39 = 80
    40 = 90
    41 = 86
    42 = 86
    43 = 84
    44 = 92
    45 = 79
    46 = 98
    47 = 92
    48 = 90
    49 = 94
}

foreach ($r in $newVals.Keys) {
    $ws.Cells.Item($r, 2).Value = $newVals[$r]
}

# --- Format B2:B49 as centered + vertically centered + wrapped ---------
# Build the combined alignment format once on a scratch cell so a single
# style (one font + one cellXf) is produced, then paint it onto the range
# via copy/paste-format (mirrors how this look was applied in Excel).
$scratch = $ws.Range("D1")
$scratch.HorizontalAlignment = -4108   # xlCenter
$scratch.VerticalAlignment = -4108     # xlCenter
$scratch.WrapText = $true
$scratch.Font.Name = "Calibri"

$scratch.Copy() | Out-Null
$dataRange = $ws.Range("B2:B49")
$dataRange.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$scratch.Clear() | Out-Null

# --- Tidy up the sheet view / selection ---------------------------------
$dataRange.Select() | Out-Null

"done"
